$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Cells are plain text in the source data, so numeric-looking Price values are
# forced to Text format before assignment (and the style is reset afterwards)
# to avoid Excel auto-converting them to floating point numbers and to keep the
# cell formatting identical to the original (no explicit style on data cells).

$ws.Range("D2").Value = "64.034.65"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.643.10"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +7.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "3.121.09"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "63.894.58"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "2.670.47"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "556.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "0.0₃0848"
$ws.Range("E33").Value = "  +5.09%  "
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.406"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0570"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0245"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0963"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.59%  "
